$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 383, shifting existing rows 383:411 down to 384:412
$ws.Rows("383:383").Insert()

# Populate the newly inserted row 383 with the new weekly record.
# Categorical columns (A,B,C,E,F,G,H,I,N,O,Q,R) repeat the values used
# throughout this "Rabanito" block; only the date/volume/price figures change.
$ws.Cells.Item(383, 1).Value = 9
$ws.Cells.Item(383, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(383, 3).Value = "Metropolitana"
$ws.Cells.Item(383, 4).Value = 45021
$ws.Cells.Item(383, 4).Style = $ws.Cells.Item(384, 4).Style
$ws.Cells.Item(383, 4).NumberFormat = $ws.Cells.Item(384, 4).NumberFormat
$ws.Cells.Item(383, 5).Value = 13
$ws.Cells.Item(383, 6).Value = 300000001
$ws.Cells.Item(383, 7).Value = "Rabanito"
$ws.Cells.Item(383, 8).Value = "Sin especificar"
$ws.Cells.Item(383, 9).Value = "Primera"
$ws.Cells.Item(383, 10).Value = 6000
$ws.Cells.Item(383, 11).Value = 3000
$ws.Cells.Item(383, 12).Value = 3000
$ws.Cells.Item(383, 13).Value = 3000
$ws.Cells.Item(383, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(383, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(383, 16).Value = 30
$ws.Cells.Item(383, 17).Value = 100
$ws.Cells.Item(383, 18).Value = "Hortaliza"
